# Auto-generated: update FFXIV Leve profit market-price figures per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4260.9287
$ws.Range("J17").Value = 3543.9167
$ws.Range("L17").Value = 10631.7501
$ws.Range("N17").Value = -10967.7501
$ws.Range("H34").Value = 3907
$ws.Range("I34").Value = 3907
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 3907
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -3704
$ws.Range("N34").ClearContents()
$ws.Range("H36").Value = 3907
$ws.Range("I36").Value = 3907
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 3907
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -3192
$ws.Range("N36").ClearContents()
$ws.Range("H40").Value = 2680
$ws.Range("I40").Value = 3000
$ws.Range("J40").Value = 2600
$ws.Range("K40").Value = 3000
$ws.Range("L40").Value = 2600
$ws.Range("M40").Value = -2825
$ws.Range("N40").Value = -2950
$ws.Range("H74").Value = 5626
$ws.Range("I74").Value = 6500
$ws.Range("J74").Value = 4752
$ws.Range("K74").Value = 6500
$ws.Range("L74").Value = 4752
$ws.Range("M74").Value = -5564
$ws.Range("N74").Value = -6624
$ws.Range("H77").Value = 5626
$ws.Range("I77").Value = 6500
$ws.Range("J77").Value = 4752
$ws.Range("K77").Value = 32500
$ws.Range("L77").Value = 23760
$ws.Range("M77").Value = -27820
$ws.Range("N77").Value = -33120
$ws.Range("H86").Value = 1337.4615
$ws.Range("I86").Value = 1282.4166
$ws.Range("K86").Value = 1282.4166
$ws.Range("M86").Value = -159.4166
$ws.Range("H89").Value = 1337.4615
$ws.Range("I89").Value = 1282.4166
$ws.Range("K89").Value = 6412.083000000001
$ws.Range("M89").Value = -796.0830000000005
$ws.Range("H98").Value = 4391.88
$ws.Range("I98").Value = 4021.4211
$ws.Range("K98").Value = 4021.4211
$ws.Range("M98").Value = -2523.4211
$ws.Range("H122").Value = 4391.88
$ws.Range("I122").Value = 4021.4211
$ws.Range("K122").Value = 12064.2633
$ws.Range("M122").Value = -9614.263300000001
$ws.Range("H131").Value = 4242.857
$ws.Range("I131").Value = 1095
$ws.Range("J131").Value = 5502
$ws.Range("K131").Value = 3285
$ws.Range("L131").Value = 16506
$ws.Range("M131").Value = 1755
$ws.Range("N131").Value = -26586
$ws.Range("H137").Value = 38125.594
$ws.Range("I137").Value = 718.9375
$ws.Range("J137").Value = 92535.27
$ws.Range("K137").Value = 2156.8125
$ws.Range("L137").Value = 277605.81
$ws.Range("M137").Value = 393.1875
$ws.Range("N137").Value = -282705.81
$ws.Range("H138").Value = 2392.699
$ws.Range("J138").Value = 2233.1965
$ws.Range("L138").Value = 6699.5895
$ws.Range("N138").Value = -16979.5895

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 86673
$ws.Range("J23").Value = 100007
$ws.Range("L23").Value = 100007
$ws.Range("N23").Value = -100525
$ws.Range("H32").Value = 4765.287
$ws.Range("I32").Value = 2611.5442
$ws.Range("K32").Value = 2611.5442
$ws.Range("M32").Value = -2324.5442
$ws.Range("H74").Value = 757.6842
$ws.Range("I74").Value = 550.9394
$ws.Range("K74").Value = 550.9394
$ws.Range("M74").Value = 323.0606
$ws.Range("H77").Value = 757.6842
$ws.Range("I77").Value = 550.9394
$ws.Range("K77").Value = 2754.697
$ws.Range("M77").Value = 1613.303
$ws.Range("H104").Value = 38799.8
$ws.Range("J104").Value = 39749.75
$ws.Range("L104").Value = 39749.75
$ws.Range("N104").Value = -46737.75
$ws.Range("H132").Value = 2071.4119
$ws.Range("I132").Value = 1626.6666
$ws.Range("J132").Value = 2571.75
$ws.Range("K132").Value = 4879.9998
$ws.Range("L132").Value = 7715.25
$ws.Range("M132").Value = -2349.9998
$ws.Range("N132").Value = -12775.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2106.2778
$ws.Range("I20").Value = 1851.25
$ws.Range("J20").Value = 2998.875
$ws.Range("K20").Value = 1851.25
$ws.Range("L20").Value = 2998.875
$ws.Range("M20").Value = -1604.25
$ws.Range("N20").Value = -3492.875
$ws.Range("H22").Value = 747.4
$ws.Range("I22").Value = 579.6667
$ws.Range("K22").Value = 579.6667
$ws.Range("M22").Value = -406.6667
$ws.Range("H94").Value = 828.3043
$ws.Range("I94").Value = 402.65
$ws.Range("K94").Value = 402.65
$ws.Range("M94").Value = 48.35000000000002
$ws.Range("H134").Value = 6855.6294
$ws.Range("I134").Value = 8102
$ws.Range("K134").Value = 24306
$ws.Range("M134").Value = -21771

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2100.9656
$ws.Range("I31").Value = 1772.6666
$ws.Range("J31").Value = 2638.182
$ws.Range("K31").Value = 1772.6666
$ws.Range("L31").Value = 2638.182
$ws.Range("M31").Value = -1477.6666
$ws.Range("N31").Value = -3228.182
$ws.Range("H34").Value = 2100.9656
$ws.Range("I34").Value = 1772.6666
$ws.Range("J34").Value = 2638.182
$ws.Range("K34").Value = 1772.6666
$ws.Range("L34").Value = 2638.182
$ws.Range("M34").Value = -1570.6666
$ws.Range("N34").Value = -3042.182
$ws.Range("H96").Value = 19000
$ws.Range("J96").Value = 19000
$ws.Range("L96").Value = 19000
$ws.Range("N96").Value = -24492

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 3466.6667
$ws.Range("I29").Value = 5100
$ws.Range("J29").Value = 200
$ws.Range("K29").Value = 15300
$ws.Range("L29").Value = 600
$ws.Range("M29").Value = -15023
$ws.Range("N29").Value = -1154
$ws.Range("H44").Value = 4501
$ws.Range("I44").Value = 5501.5
$ws.Range("K44").Value = 16504.5
$ws.Range("M44").Value = -16106.5
$ws.Range("H121").Value = 668.2778
$ws.Range("I121").Value = 494.8
$ws.Range("J121").Value = 735
$ws.Range("K121").Value = 1484.4
$ws.Range("L121").Value = 2205
$ws.Range("M121").Value = -174.4000000000001
$ws.Range("N121").Value = -4825
$ws.Range("H132").Value = 10940.4
$ws.Range("I132").Value = 984
$ws.Range("J132").Value = 25875
$ws.Range("K132").Value = 8856
$ws.Range("L132").Value = 232875
$ws.Range("M132").Value = -6326
$ws.Range("N132").Value = -237935
$ws.Range("H140").Value = 2485.389
$ws.Range("I140").Value = 1877.4375
$ws.Range("K140").Value = 5632.3125
$ws.Range("M140").Value = -452.3125
$ws.Range("H141").Value = 2890.5
$ws.Range("I141").Value = 2816.95
$ws.Range("J141").Value = 3258.25
$ws.Range("K141").Value = 8450.849999999999
$ws.Range("L141").Value = 9774.75
$ws.Range("M141").Value = -3270.849999999999
$ws.Range("N141").Value = -20134.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2830975.8
$ws.Range("I126").Value = 3089778.5
$ws.Range("K126").Value = 9269335.5
$ws.Range("M126").Value = -9266865.5
$ws.Range("H132").Value = 2140425.5
$ws.Range("I132").Value = 2750381
$ws.Range("J132").Value = 5581.5
$ws.Range("K132").Value = 8251143
$ws.Range("L132").Value = 16744.5
$ws.Range("M132").Value = -8248613
$ws.Range("N132").Value = -21804.5
$ws.Range("H136").Value = 5020.826
$ws.Range("J136").Value = 5020.826
$ws.Range("L136").Value = 15062.478
$ws.Range("N136").Value = -20162.478

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2983.8572
$ws.Range("I61").Value = 2972.5
$ws.Range("K61").Value = 2972.5
$ws.Range("M61").Value = -2770.5
$ws.Range("H113").Value = 2983.8572
$ws.Range("I113").Value = 2972.5
$ws.Range("K113").Value = 2972.5
$ws.Range("M113").Value = -802.5
$ws.Range("H127").Value = 39755.668
$ws.Range("J127").Value = 39755.668
$ws.Range("L127").Value = 39755.668
$ws.Range("N127").Value = -49675.668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 13980
$ws.Range("J31").Value = 13980
$ws.Range("L31").Value = 13980
$ws.Range("N31").Value = -14676
